# Rebuilt with cut corner
# Update the parameter sheet: change several values, rename some
# parameters (removing MPPC/mounting-screw rows and adding new
# corner-chamfer / mounting-spacing / silvering / rubberized-coating
# rows), and delete the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: MPPCClearence -> FiberEdgeOffsett -----------------------
$ws.Range("A10").Value = "FiberEdgeOffsett"
$ws.Range("B10").Value = 5

# --- Row 11: MPPCSensorClearence -> LoopLargeTrackRatio --------------
$ws.Range("A11").Value = "LoopLargeTrackRatio"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "ul"

# --- Row 12: MPPCSensorOffsett -> LoopSmallTrackDiameter -------------
$ws.Range("A12").Value = "LoopSmallTrackDiameter"
$ws.Range("B12").Value = 100

# --- Row 13: MPPCWidth -> FiberSensorClerence -------------------------
$ws.Range("A13").Value = "FiberSensorClerence"
$ws.Range("B13").Value = 0.5

# --- Row 15: FiberEdgeOffsett -> cornerChamfer -----------------------
$ws.Range("A15").Value = "cornerChamfer"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "in"

# --- Row 14: MPPCHeight -> mountingSpacing ---------------------------
$ws.Range("A14").Value = "mountingSpacing"
$ws.Range("B14").Value = 18

# --- Row 16: MPPCDepth -> silveringThickness -------------------------
$ws.Range("A16").Value = "silveringThickness"
$ws.Range("B16").Value = 0.2

# --- Row 17: LongFiberAngle -> rubberizedCoating ---------------------
$ws.Range("A17").Value = "rubberizedCoating"
$ws.Range("B17").Value = 20
$ws.Range("C17").Value = "mil"

# --- Remove now-unused rows 18-24 ------------------------------------
$ws.Range("A18:C24").ClearContents()

# --- Update selection to match the authored commit -------------------
$ws.Range("B18").Select()
